$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 0.1807692307692308
    "C2" = 0.5692307692307692
    "P2" = 0.1346153846153846
    "S2" = 0.1153846153846154
    "C3" = 0.02580645161290323
    "J3" = 0.01290322580645161
    "P3" = 0.7935483870967742
    "S3" = 0.167741935483871
    "J4" = 0.02439024390243903
    "P4" = 0.6585365853658537
    "S4" = 0.3170731707317073
    "B6" = 0.04694835680751173
    "D6" = 0.01408450704225352
    "F6" = 0.06572769953051644
    "J6" = 0.2723004694835681
    "O6" = 0.01408450704225352
    "Q6" = 0.2065727699530517
    "R6" = 0.07981220657276995
    "S6" = 0.3004694835680751
    "B7" = 0.09045226130653267
    "D7" = 0.02010050251256281
    "F7" = 0.05527638190954774
    "J7" = 0.07537688442211055
    "O7" = 0.02010050251256281
    "Q7" = 0.1608040201005025
    "R7" = 0.07035175879396985
    "S7" = 0.507537688442211
    "B8" = 0.07048458149779736
    "D8" = 0.01762114537444934
    "F8" = 0.07048458149779736
    "J8" = 0.1475770925110132
    "O8" = 0.05286343612334802
    "Q8" = 0.1850220264317181
    "R8" = 0.07929515418502203
    "S8" = 0.3766519823788546
    "B9" = 0.08695652173913043
    "D9" = 0.004830917874396135
    "F9" = 0.05314009661835749
    "J9" = 0.1497584541062802
    "O9" = 0.04347826086956522
    "Q9" = 0.2270531400966184
    "R9" = 0.0821256038647343
    "S9" = 0.3526570048309179
    "B10" = 0.1044176706827309
    "D10" = 0.02168674698795181
    "E10" = 0.0008032128514056225
    "F10" = 0.06506024096385542
    "J10" = 0.1132530120481928
    "O10" = 0.02329317269076305
    "Q10" = 0.236144578313253
    "R10" = 0.09718875502008033
    "S10" = 0.3381526104417671
    "G11" = 0.1254237288135593
    "J11" = 0.09830508474576272
    "K11" = 0.1830508474576271
    "L11" = 0.576271186440678
    "S11" = 0.01694915254237288
    "F12" = 0.005681818181818182
    "G12" = 0.7102272727272727
    "J12" = 0.1931818181818182
    "L12" = 0.02840909090909091
    "S12" = 0.0625
    "G13" = 0.6065573770491803
    "J13" = 0.3442622950819672
    "S13" = 0.04918032786885246
    "F15" = 0.01298701298701299
    "H15" = 0.1255411255411255
    "I15" = 0.06060606060606061
    "J15" = 0.2640692640692641
    "K15" = 0.06493506493506493
    "M15" = 0.01298701298701299
    "N15" = 0.008658008658008658
    "O15" = 0.06926406926406926
    "S15" = 0.3809523809523809
    "F16" = 0.01104972375690608
    "H16" = 0.1602209944751381
    "I16" = 0.09392265193370165
    "J16" = 0.3535911602209945
    "K16" = 0.1160220994475138
    "M16" = 0.06629834254143646
    "O16" = 0.04419889502762431
    "S16" = 0.1546961325966851
    "F17" = 0.01996007984031936
    "H17" = 0.1596806387225549
    "I17" = 0.1017964071856287
    "J17" = 0.4251497005988024
    "K17" = 0.08383233532934131
    "M17" = 0.01796407185628742
    "O17" = 0.03992015968063872
    "S17" = 0.1516966067864272
    "F18" = 0.009708737864077669
    "H18" = 0.2233009708737864
    "I18" = 0.0970873786407767
    "J18" = 0.3786407766990291
    "K18" = 0.07766990291262135
    "M18" = 0.02912621359223301
    "N18" = 0.004854368932038835
    "O18" = 0.02912621359223301
    "S18" = 0.1504854368932039
    "F19" = 0.01557632398753894
    "H19" = 0.2110591900311526
    "I19" = 0.08333333333333333
    "J19" = 0.3543613707165109
    "K19" = 0.1129283489096573
    "M19" = 0.0264797507788162
    "N19" = 0.000778816199376947
    "O19" = 0.06619937694704049
    "S19" = 0.1292834890965732
}

foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}
